$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195, shifting existing rows 195:290 down to 196:291
$ws.Rows("195:195").Insert()

# Populate new row 195 with values. Columns A,B,C,E,F,G,H,I,R are constant
# across this data block (same market / product), matching row 196 (the row
# that used to be row 195 before the insert).
$ws.Range("A195").Value = 4
$ws.Range("B195").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C195").Value = "Los Lagos"
$ws.Range("D195").Value = 44806
$ws.Range("E195").Value = 10
$ws.Range("F195").Value = 100112032
$ws.Range("G195").Value = "Zapallo italiano"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 200
$ws.Range("K195").Value = 28000
$ws.Range("L195").Value = 30000
$ws.Range("M195").Value = 29000
$ws.Range("N195").Value = '$/caja 50 unidades'
$ws.Range("O195").Value = "Región de Arica y Parinacota"
$ws.Range("P195").Value = 580
$ws.Range("Q195").Value = 50
$ws.Range("R195").Value = "Hortaliza"

# Apply the same date number format used by the other rows in column D
$ws.Range("D195").NumberFormat = $ws.Range("D196").NumberFormat
